$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 1.468507333333333
$ws.Range("H2").Value = 4.405521999999999
$ws.Range("I2").Value = 0.005118279455112885
$ws.Range("J2").Value = 0.005118279455112885
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 1.101036666666667
$ws.Range("N2").Value = 3.30311
$ws.Range("O2").Value = 0.007197401111328883
$ws.Range("P2").Value = 0.007197401111328884
$ws.Range("Q2").Value = 1.616880419268889
$ws.Range("R2").Value = 14.55192377342
$ws.Range("S2").Value = 0.00003683831023832127
$ws.Range("T2").Value = 0.00003683831023832127

$ws.Range("G3").Value = 1.468507333333333
$ws.Range("H3").Value = 4.405521999999999
$ws.Range("I3").Value = 0.005118279455112885
$ws.Range("J3").Value = 0.005118279455112885
$ws.Range("O3").Value = 0.0007216619689517899
$ws.Range("P3").Value = 0.00072166196895179
$ws.Range("Q3").Value = 0.1621197830828889
$ws.Range("R3").Value = 1.459078047746
$ws.Range("S3").Value = 0.000003693667629222259
$ws.Range("T3").Value = 0.000003693667629222259

$ws.Range("G4").Value = 1.468507333333333
$ws.Range("H4").Value = 4.405521999999999
$ws.Range("I4").Value = 0.005118279455112885
$ws.Range("J4").Value = 0.005118279455112885
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.1677316666666666
$ws.Range("N4").Value = 0.5031949999999999
$ws.Range("O4").Value = 0.00109645039136303
$ws.Range("P4").Value = 0.00109645039136303
$ws.Range("Q4").Value = 0.2463151825322221
$ws.Range("R4").Value = 2.21683664279
$ws.Range("S4").Value = 0.000005611939511663876
$ws.Range("T4").Value = 0.000005611939511663878

$ws.Range("G5").Value = 1.468507333333333
$ws.Range("H5").Value = 4.405521999999999
$ws.Range("I5").Value = 0.005118279455112885
$ws.Range("J5").Value = 0.005118279455112885
$ws.Range("M5").Value = 151.5978113333333
$ws.Range("N5").Value = 454.793434
$ws.Range("O5").Value = 0.9909844865283564
$ws.Range("P5").Value = 0.9909844865283564
$ws.Range("Q5").Value = 222.6224976602831
$ws.Range("R5").Value = 2003.602478942548
$ws.Range("S5").Value = 0.005072135537733678
$ws.Range("T5").Value = 0.005072135537733678

$ws.Range("I6").Value = 0.9046276674881553
$ws.Range("J6").Value = 0.9046276674881553
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 1.101036666666667
$ws.Range("N6").Value = 3.30311
$ws.Range("O6").Value = 0.007197401111328883
$ws.Range("P6").Value = 0.007197401111328884
$ws.Range("Q6").Value = 285.7746973603312
$ws.Range("R6").Value = 2571.972276242981
$ws.Range("S6").Value = 0.006510968179318104
$ws.Range("T6").Value = 0.006510968179318104

$ws.Range("I7").Value = 0.9046276674881553
$ws.Range("J7").Value = 0.9046276674881553
$ws.Range("O7").Value = 0.0007216619689517899
$ws.Range("P7").Value = 0.00072166196895179
$ws.Range("Q7").Value = 28.65377760439712
$ws.Range("S7").Value = 0.0006528353836877672
$ws.Range("T7").Value = 0.0006528353836877673

$ws.Range("I8").Value = 0.9046276674881553
$ws.Range("J8").Value = 0.9046276674881553
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 0.1677316666666666
$ws.Range("N8").Value = 0.5031949999999999
$ws.Range("O8").Value = 0.00109645039136303
$ws.Range("P8").Value = 0.00109645039136303
$ws.Range("Q8").Value = 43.53485013766778
$ws.Range("R8").Value = 391.81365123901
$ws.Range("S8").Value = 0.0009918793600552124
$ws.Range("T8").Value = 0.0009918793600552126

$ws.Range("I9").Value = 0.9046276674881553
$ws.Range("J9").Value = 0.9046276674881553
$ws.Range("M9").Value = 151.5978113333333
$ws.Range("N9").Value = 454.793434
$ws.Range("O9").Value = 0.9909844865283564
$ws.Range("P9").Value = 0.9909844865283564
$ws.Range("Q9").Value = 39347.29874657997
$ws.Range("R9").Value = 354125.6887192196
$ws.Range("S9").Value = 0.8964719845650942
$ws.Range("T9").Value = 0.8964719845650942

$ws.Range("G10").Value = 0.5890733333333333
$ws.Range("H10").Value = 1.76722
$ws.Range("I10").Value = 0.002053133730501083
$ws.Range("J10").Value = 0.002053133730501083
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 1.101036666666667
$ws.Range("N10").Value = 3.30311
$ws.Range("O10").Value = 0.007197401111328883
$ws.Range("P10").Value = 0.007197401111328884
$ws.Range("Q10").Value = 0.6485913393555556
$ws.Range("R10").Value = 5.8373220542
$ws.Range("S10").Value = 0.00001477722699361531
$ws.Range("T10").Value = 0.00001477722699361531

$ws.Range("G11").Value = 0.5890733333333333
$ws.Range("H11").Value = 1.76722
$ws.Range("I11").Value = 0.002053133730501083
$ws.Range("J11").Value = 0.002053133730501083
$ws.Range("O11").Value = 0.0007216619689517899
$ws.Range("P11").Value = 0.00072166196895179
$ws.Range("Q11").Value = 0.06503232149555556
$ws.Range("R11").Value = 0.5852908934600001
$ws.Range("S11").Value = 0.000001481668530474745
$ws.Range("T11").Value = 0.000001481668530474745

$ws.Range("G12").Value = 0.5890733333333333
$ws.Range("H12").Value = 1.76722
$ws.Range("I12").Value = 0.002053133730501083
$ws.Range("J12").Value = 0.002053133730501083
$ws.Range("K12").Value = 2
$ws.Range("L12").Value = 0.6666666666666666
$ws.Range("M12").Value = 0.1677316666666666
$ws.Range("N12").Value = 0.5031949999999999
$ws.Range("O12").Value = 0.00109645039136303
$ws.Range("P12").Value = 0.00109645039136303
$ws.Range("Q12").Value = 0.09880625198888887
$ws.Range("R12").Value = 0.8892562678999999
$ws.Range("S12").Value = 0.00000225115928232855
$ws.Range("T12").Value = 0.00000225115928232855

$ws.Range("G13").Value = 0.5890733333333333
$ws.Range("H13").Value = 1.76722
$ws.Range("I13").Value = 0.002053133730501083
$ws.Range("J13").Value = 0.002053133730501083
$ws.Range("M13").Value = 151.5978113333333
$ws.Range("N13").Value = 454.793434
$ws.Range("O13").Value = 0.9909844865283564
$ws.Range("P13").Value = 0.9909844865283564
$ws.Range("Q13").Value = 89.30222804816445
$ws.Range("R13").Value = 803.72005243348
$ws.Range("S13").Value = 0.002034623675694665
$ws.Range("T13").Value = 0.002034623675694665

$ws.Range("G14").Value = 25.306101
$ws.Range("H14").Value = 75.91830299999999
$ws.Range("I14").Value = 0.0882009193262308
$ws.Range("J14").Value = 0.0882009193262308
$ws.Range("K14").Value = 3
$ws.Range("L14").Value = 1
$ws.Range("M14").Value = 1.101036666666667
$ws.Range("N14").Value = 3.30311
$ws.Range("O14").Value = 0.007197401111328883
$ws.Range("P14").Value = 0.007197401111328884
$ws.Range("Q14").Value = 27.86294509137
$ws.Range("R14").Value = 250.76650582233
$ws.Range("S14").Value = 0.0006348173947788426
$ws.Range("T14").Value = 0.0006348173947788427

$ws.Range("G15").Value = 25.306101
$ws.Range("H15").Value = 75.91830299999999
$ws.Range("I15").Value = 0.0882009193262308
$ws.Range("J15").Value = 0.0882009193262308
$ws.Range("O15").Value = 0.0007216619689517899
$ws.Range("P15").Value = 0.00072166196895179
$ws.Range("Q15").Value = 2.793734502831
$ws.Range("R15").Value = 25.143610525479
$ws.Range("S15").Value = 0.00006365124910432569
$ws.Range("T15").Value = 0.0000636512491043257

$ws.Range("G16").Value = 25.306101
$ws.Range("H16").Value = 75.91830299999999
$ws.Range("I16").Value = 0.0882009193262308
$ws.Range("J16").Value = 0.0882009193262308
$ws.Range("K16").Value = 2
$ws.Range("L16").Value = 0.6666666666666666
$ws.Range("M16").Value = 0.1677316666666666
$ws.Range("N16").Value = 0.5031949999999999
$ws.Range("O16").Value = 0.00109645039136303
$ws.Range("P16").Value = 0.00109645039136303
$ws.Range("Q16").Value = 4.244634497564999
$ws.Range("R16").Value = 38.20171047808499
$ws.Range("S16").Value = 0.00009670793251382475
$ws.Range("T16").Value = 0.00009670793251382477

$ws.Range("G17").Value = 25.306101
$ws.Range("H17").Value = 75.91830299999999
$ws.Range("I17").Value = 0.0882009193262308
$ws.Range("J17").Value = 0.0882009193262308
$ws.Range("M17").Value = 151.5978113333333
$ws.Range("N17").Value = 454.793434
$ws.Range("O17").Value = 0.9909844865283564
$ws.Range("P17").Value = 0.9909844865283564
$ws.Range("Q17").Value = 3836.3495249802777
$ws.Range("R17").Value = 34527.1457248225
$ws.Range("S17").Value = 0.0874057427498338
$ws.Range("T17").Value = 0.0874057427498338

Write-Output "done"